# "2 and 11 are fixed" - correct the answer-key table (variants 1-3) for
# a handful of question cells that had the wrong letter.
$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

function Set-AnswerCell($row, $col, $letter) {
    $cell = $table.Cell($row, $col)
    $cell.Range.Text = $letter
}

# Variant 1 (row 2 in the table, row 1 has the header)
Set-AnswerCell 2 2 "А"
Set-AnswerCell 2 4 "В"
Set-AnswerCell 2 6 "А"

# Variant 2 (row 3)
Set-AnswerCell 3 2 "Г"
Set-AnswerCell 3 4 "Б"
Set-AnswerCell 3 5 "В"
Set-AnswerCell 3 6 "В"
Set-AnswerCell 3 7 "Г"

# Variant 3 (row 4)
Set-AnswerCell 4 3 "А"
Set-AnswerCell 4 4 "В"
Set-AnswerCell 4 5 "В"
Set-AnswerCell 4 6 "В"
